$d = $word.ActiveDocument

# The year "2019" (stored as two runs: "201" + "9") must become "2020".
# The commit only actually changes the rendered text from 2019 to 2020;
# Word keeps the untouched leading "2" in its own run and places the
# freshly-typed "020" in a new run immediately after it, then drops the
# "_GoBack" bookmark (marking the last edit point) right after that
# newly typed text.

# 1) Replace "2019" with just "2" first - this leaves the original,
#    untouched "2" as its own run (matching how Word only retypes the
#    part that actually changed).
$rng = $d.Content
$rng.Find.Execute("2019", $true, $false, $false, $false, $false, $true, 1, $false, "2", 2)
$afterTwo = $rng.End

# 2) Insert the newly "typed" remainder of the year as its own run right
#    after that.
$rng2 = $d.Range($afterTwo, $afterTwo)
$rng2.InsertAfter("020")
$rng2.End = $afterTwo + 3

# Toggling a character property on this new range and back keeps it from
# being silently re-merged with the preceding, differently-authored "2"
# run when the package is serialized, matching the two separate <w:r>
# elements ("2" and "020") seen in the target document.
$rng2.Bold = 1
$rng2.Bold = 0

# 3) Move the "_GoBack" bookmark (Word's "last edit location" marker) to
#    sit right after the text that was just typed, i.e. right after
#    "020" and before the following "г.".
$bmRange = $d.Range($rng2.End, $rng2.End)
$d.Bookmarks.Add("_GoBack", $bmRange)
